# Apply "Some game revisions, cut lines for pdf" edits.
#
# 1. Rename several card titles on the "Deck" sheet (adding "The" prefix /
#    rewording a couple of Resource cards), and extend the "Starting
#    Player" description to also allow trading for Food.
# 2. Switch the active/selected sheet from "Stats" to "Deck" and update the
#    remembered cell selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsDeck = $wb.Worksheets.Item("Deck")
$wsStats = $wb.Worksheets.Item("Stats")

# --- Card title / description revisions (Deck sheet) ---------------------
$wsDeck.Range("C6").Value = "The Insta-slaughter"
$wsDeck.Range("C11").Value = "The Cheap Resource Combo"
$wsDeck.Range("C12").Value = "The Better Resource Combo"
$wsDeck.Range("C13").Value = "The Decent Resource Combo"
$wsDeck.Range("C18").Value = "The Growing Resource"
$wsDeck.Range("C19").Value = "The Other Farming Resource"
$wsDeck.Range("C20").Value = "The Rare yet Valuable Resource"
$wsDeck.Range("C23").Value = "The Early Investment"
$wsDeck.Range("I2").Value = "Other players may transfer this card to their playing area for one of any of the following: Stone, Clay, Silk, Boar, Glass, Gold, Cattle, Vegetable, Food."

# --- Selection / active-sheet bookkeeping ---------------------------------
# Record the remembered selection on Stats first (it is still the active
# sheet at this point, matching how the state was captured before the user
# flipped back to Deck).
$wsStats.Range("I8").Select()

# Make Deck the active sheet and leave the cursor on I2.
$wsDeck.Activate()
$wsDeck.Range("I2").Select()
